$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '29.377.87'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -3.30%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '1.986.71'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -5.28%  '

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '1.019'
$c.Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +1.68%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '328.88'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -4.24%  '

$ws.Cells.Item(6, 5).Value = '  +1.37%  '

$ws.Cells.Item(7, 5).Value = '  -6.30%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.4219'
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -4.96%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '53.97'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -1.48%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.08883'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -5.27%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '1.106'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -5.51%  '

$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '2.145.23'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -0.51%  '

$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '23.09'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -6.71%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '7.900'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -7.58%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '6.437'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -6.81%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '1.018'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +1.56%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '93.91'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -7.67%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '0.00001105'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -4.81%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '0.06746'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +0.83%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '19.30'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -8.86%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '1.016'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +1.39%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '5.946'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -6.04%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '29.462.77'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -3.16%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '11.96'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -4.65%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '2.316'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.17%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '20.68'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -5.56%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '156.31'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -3.99%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '6.207'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -8.31%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.292'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -9.10%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '127.10'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.84%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.050'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -8.25%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '0.09916'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -6.21%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '1.514'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -8.64%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '5.800'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -7.30%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '3.793'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -1.60%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.02442'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -7.30%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '9.196'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -9.79%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '0.06365'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -6.38%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '1.289'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -4.45%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '0.6508'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -7.38%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '11.56'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -8.71%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.2027'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -8.72%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '1.015'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.45%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.6268'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -8.71%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '13.41'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -6.76%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '2.196'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -6.16%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '1.289'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -7.35%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '3.497'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -3.88%  '

$ws.Cells.Item(49, 5).Value = '  -1.77%  '

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '0.06936'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -4.27%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '1.124'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -8.81%  '
